# Correção no teste t para uma população - valor_p
# Atualiza os valores de A2:A248 com a nova amostra e remove as linhas A249:A306
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    9.6,
    13.61,
    10.5,
    5.63,
    10.32,
    14.17,
    14.27,
    9.359999999999999,
    1.18,
    14.69,
    14.69,
    12.31,
    9.6,
    11.75,
    7.24,
    9.800000000000001,
    3.58,
    9.800000000000001,
    5.08,
    13.07,
    4.9,
    4.31,
    3.67,
    15.66,
    7.69,
    5.08,
    4.85,
    10.89,
    12.45,
    9.56,
    4.98,
    0.86,
    24.28,
    8.859999999999999,
    6.66,
    13.44,
    15.57,
    9.43,
    7.13,
    8.5,
    13.72,
    12.63,
    12.58,
    16.22,
    13.69,
    12.15,
    11.89,
    5.92,
    15.18,
    11.56,
    6.73,
    10.94,
    7.85,
    14.09,
    9.74,
    12.43,
    15.44,
    14.98,
    5.63,
    18.72,
    11.61,
    14.72,
    10.63,
    16.2,
    11.74,
    7.41,
    11.8,
    1,
    3.14,
    13.3,
    15.71,
    8.529999999999999,
    7.34,
    0.98,
    15.29,
    8.949999999999999,
    13.98,
    11.43,
    11.75,
    5.58,
    15.43,
    8.35,
    9.970000000000001,
    10.63,
    16.26,
    21.29,
    10.26,
    14.47,
    6.05,
    7.85,
    12.66,
    13.91,
    19.43,
    -2.68,
    17.52,
    8.470000000000001,
    9.279999999999999,
    7.28,
    7.19,
    2.94,
    7.86,
    12.98,
    6.47,
    4.52,
    7.38,
    9.34,
    12.31,
    7.44,
    7.82,
    16.82,
    8.25,
    13.66,
    14.01,
    1.55,
    13.2,
    11.36,
    7.59,
    14.44,
    13.34,
    9.94,
    7.27,
    9.98,
    14.58,
    18.05,
    12.12,
    13.92,
    9.640000000000001,
    7.7,
    8.98,
    11.22,
    6.6,
    5.14,
    14.97,
    12.82,
    9.470000000000001,
    11.93,
    17.85,
    17.33,
    16.19,
    2.61,
    10.89,
    15.72,
    11.85,
    16.14,
    12.15,
    -1.05,
    7.69,
    14.21,
    -0.33,
    14.86,
    7.1,
    4.6,
    7.06,
    14.48,
    13.79,
    5.9,
    13.31,
    10.09,
    10.68,
    7.51,
    7.63,
    11.88,
    16.58,
    0.85,
    14.03,
    7.88,
    15.41,
    12.05,
    13.74,
    4.29,
    9.460000000000001,
    16.16,
    10.94,
    13.4,
    7.01,
    7.16,
    11.46,
    15.46,
    8.85,
    10.32,
    8.57,
    8.92,
    10.74,
    12.71,
    -1.89,
    8.449999999999999,
    15.5,
    12.92,
    10.69,
    13.1,
    0.28,
    17.96,
    18.5,
    13.23,
    14.55,
    7.41,
    7.04,
    8.81,
    13.49,
    14.93,
    9.49,
    1.6,
    7.2,
    6.62,
    8.609999999999999,
    7.17,
    10.31,
    7.12,
    6.91,
    9.529999999999999,
    4.05,
    6.05,
    10.26,
    18.54,
    17.85,
    9.67,
    3.51,
    11.03,
    18.17,
    5.6,
    8.720000000000001,
    15.98,
    8.73,
    10.04,
    17.17,
    10.26,
    7.14,
    18.3,
    3.24,
    10.58,
    17.06,
    4.65,
    9.050000000000001,
    14.81,
    4.84,
    6.05,
    10.31,
    11.04,
    13.62,
    9.42,
    5.63,
    12.82,
    10.73,
    13.01,
    12.67,
    3,
    11.35
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $values[$i]
}

# Remove as linhas que sobraram (249:306) para que a dimensao fique A1:A248
$ws.Rows("249:306").Delete()
